$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes -----------------------------------------------------
# Update the author list text (in place - same cell, text amended)
$ws.Range("H2").Value = "Daniela Subotic, Noémi Villars-Amberg"

# Add new column I: "Authorship Resource" holding the same author text
$ws.Range("I1").Value = "Authorship Resource"
$ws.Range("I2").Value = "Daniela Subotic, Noémi Villars-Amberg"

# Give the new header cell (I1) the same formatting as the other header
# cells (bold font + bottom border, left aligned) by copying the format
# from H1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Make the new column a reasonable width
$ws.Columns("I").ColumnWidth = 44.6667

# --- Style clean-up -----------------------------------------------------
# The data row cells had an explicit "general alignment" style; reset them
# back to the workbook default (Normal) style, keeping the D2 header-style
# cell untouched.
$ws.Range("A2:C2").Style = "Normal"
$ws.Range("E2:H2").Style = "Normal"

# Restore the cursor position like in the authored workbook
$ws.Range("C15").Select() | Out-Null
